# Update alternate-cost / alternate-flag tables and refresh sheet selections.
$wb = $excel.ActiveWorkbook

# --- "AternateCosts" sheet: move the cost of 3 from P1->P2b to P2b->P2a ---
$wsCosts = $wb.Worksheets.Item("AternateCosts")
$wsCosts.Range("D2").Value = 0
$wsCosts.Range("C4").Value = 3
$wsCosts.Select()
$wsCosts.Range("D5").Select()

# --- "Alternates" sheet: move the flag of 1 from P2a->P2b to P2b->P2a ---
$wsAlt = $wb.Worksheets.Item("Alternates")
$wsAlt.Range("D3").Value = 0
$wsAlt.Range("C4").Value = 1
$wsAlt.Select()
$wsAlt.Range("E5").Select()
